# Auto-generated Excel COM-interop script
# Updates numeric cell values in sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
# to reflect refreshed market-price data (chore: update Sheets via scheduled runner)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 45996.668
$ws.Range("J93").Value = 45996.668
$ws.Range("L93").Value = 45996.668
$ws.Range("N93").Value = -50988.668

$ws.Range("H137").Value = 3669150.2
$ws.Range("I137").Value = 8548128
$ws.Range("J137").Value = 9916.833000000001
$ws.Range("K137").Value = 25644384
$ws.Range("L137").Value = 29750.499
$ws.Range("M137").Value = -25641834
$ws.Range("N137").Value = -34850.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7728.4287
$ws.Range("I32").Value = 3857
$ws.Range("J32").Value = 30957
$ws.Range("K32").Value = 3857
$ws.Range("L32").Value = 30957
$ws.Range("M32").Value = -3570
$ws.Range("N32").Value = -31531

$ws.Range("H44").Value = 37630
$ws.Range("J44").Value = 37630
$ws.Range("L44").Value = 37630
$ws.Range("N44").Value = -38606

$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H61").Value = 2471.7354
$ws.Range("I61").Value = 1466.25
$ws.Range("K61").Value = 1466.25
$ws.Range("M61").Value = -1254.25

$ws.Range("H74").Value = 1556.0392
$ws.Range("J74").Value = 2630.3333
$ws.Range("L74").Value = 2630.3333
$ws.Range("N74").Value = -4378.3333

$ws.Range("H77").Value = 1556.0392
$ws.Range("J77").Value = 2630.3333
$ws.Range("L77").Value = 13151.6665
$ws.Range("N77").Value = -21887.6665

$ws.Range("H80").Value = 47500
$ws.Range("J80").Value = 47500
$ws.Range("L80").Value = 47500
$ws.Range("N80").Value = -49496

$ws.Range("H83").Value = 47500
$ws.Range("J83").Value = 47500
$ws.Range("L83").Value = 142500
$ws.Range("N83").Value = -152484

$ws.Range("H123").Value = 40429
$ws.Range("J123").Value = 40429
$ws.Range("L123").Value = 40429
$ws.Range("N123").Value = -50229

$ws.Range("H132").Value = 19233756
$ws.Range("I132").Value = 38463716
$ws.Range("J132").Value = 3798.1538
$ws.Range("K132").Value = 115391148
$ws.Range("L132").Value = 11394.4614
$ws.Range("M132").Value = -115388618
$ws.Range("N132").Value = -16454.4614

$ws.Range("H136").Value = 2471.7354
$ws.Range("I136").Value = 1466.25
$ws.Range("K136").Value = 4398.75
$ws.Range("M136").Value = -1848.75

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H141").Value = 39387.6
$ws.Range("J141").Value = 39387.6
$ws.Range("L141").Value = 39387.6
$ws.Range("N141").Value = -49747.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 12103.417
$ws.Range("I75").Value = 4206.8335
$ws.Range("K75").Value = 4206.8335
$ws.Range("M75").Value = -3270.8335

$ws.Range("H78").Value = 12103.417
$ws.Range("I78").Value = 4206.8335
$ws.Range("K78").Value = 12620.5005
$ws.Range("M78").Value = -7940.500499999998

$ws.Range("H134").Value = 2594.2856
$ws.Range("I134").Value = 2086.9048
$ws.Range("K134").Value = 6260.714399999999
$ws.Range("M134").Value = -3725.714399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49999
$ws.Range("J20").Value = 49999
$ws.Range("L20").Value = 49999
$ws.Range("N20").Value = -50471

$ws.Range("H30").Value = 49999
$ws.Range("J30").Value = 49999
$ws.Range("L30").Value = 49999
$ws.Range("N30").Value = -50181

$ws.Range("H31").Value = 5381570.5
$ws.Range("I31").Value = 1646.0741
$ws.Range("J31").Value = 9531798
$ws.Range("K31").Value = 1646.0741
$ws.Range("L31").Value = 9531798
$ws.Range("M31").Value = -1351.0741
$ws.Range("N31").Value = -9532388

$ws.Range("H34").Value = 5381570.5
$ws.Range("I34").Value = 1646.0741
$ws.Range("J34").Value = 9531798
$ws.Range("K34").Value = 1646.0741
$ws.Range("L34").Value = 9531798
$ws.Range("M34").Value = -1444.0741
$ws.Range("N34").Value = -9532202

$ws.Range("H128").Value = 49999
$ws.Range("J128").Value = 49999
$ws.Range("L128").Value = 49999
$ws.Range("N128").Value = -59959

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 4241.8623
$ws.Range("J122").Value = 8245.286
$ws.Range("L122").Value = 74207.57399999999
$ws.Range("N122").Value = -79107.57399999999

$ws.Range("H131").Value = 2394.2104
$ws.Range("I131").Value = 561.25
$ws.Range("J131").Value = 3727.2727
$ws.Range("K131").Value = 1683.75
$ws.Range("L131").Value = 11181.8181
$ws.Range("M131").Value = 3356.25
$ws.Range("N131").Value = -21261.8181

$ws.Range("H139").Value = 117746.59
$ws.Range("I139").Value = 241088.08
$ws.Range("J139").Value = 3215.2144
$ws.Range("K139").Value = 723264.24
$ws.Range("L139").Value = 9645.643199999999
$ws.Range("M139").Value = -718124.24
$ws.Range("N139").Value = -19925.6432

$ws.Range("H141").Value = 71645360
$ws.Range("I141").Value = 91183630
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 273550890
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -273545710
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 840.8333
$ws.Range("J102").Value = 866.6667
$ws.Range("L102").Value = 866.6667
$ws.Range("N102").Value = -4110.6667

$ws.Range("H132").Value = 38464440
$ws.Range("I132").Value = 111113420
$ws.Range("J132").Value = 3215.1177
$ws.Range("K132").Value = 333340260
$ws.Range("L132").Value = 9645.3531
$ws.Range("M132").Value = -333337730
$ws.Range("N132").Value = -14705.3531

$ws.Range("H133").Value = 61554.668
$ws.Range("J133").Value = 61554.668
$ws.Range("L133").Value = 61554.668
$ws.Range("N133").Value = -71674.66800000001

$ws.Range("H135").Value = 76330
$ws.Range("J135").Value = 76330
$ws.Range("L135").Value = 76330
$ws.Range("N135").Value = -86470

$ws.Range("H140").Value = 40515.8
$ws.Range("J140").Value = 40515.8
$ws.Range("L140").Value = 40515.8
$ws.Range("N140").Value = -50875.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3663.9565
$ws.Range("I68").Value = 3592.8235
$ws.Range("J68").Value = 3865.5
$ws.Range("K68").Value = 3592.8235
$ws.Range("L68").Value = 3865.5
$ws.Range("M68").Value = -2843.8235
$ws.Range("N68").Value = -5363.5

$ws.Range("H71").Value = 3663.9565
$ws.Range("I71").Value = 3592.8235
$ws.Range("J71").Value = 3865.5
$ws.Range("K71").Value = 17964.1175
$ws.Range("L71").Value = 19327.5
$ws.Range("M71").Value = -14220.1175
$ws.Range("N71").Value = -26815.5

$ws.Range("H132").Value = 4072.8
$ws.Range("I132").Value = 2692.9375
$ws.Range("J132").Value = 5649.7856
$ws.Range("K132").Value = 8078.8125
$ws.Range("L132").Value = 16949.3568
$ws.Range("M132").Value = -5548.8125
$ws.Range("N132").Value = -22009.3568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2788.889
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2788.889
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2788.889
$ws.Range("N62").Value = -4036.889
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 2788.889
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2788.889
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 13944.445
$ws.Range("N65").Value = -20184.445
$ws.Range("M65").ClearContents()

$ws.Range("H107").Value = 5556860
$ws.Range("I107").Value = 1113.7693
$ws.Range("J107").Value = 20001800
$ws.Range("K107").Value = 3341.3079
$ws.Range("L107").Value = 60005400
$ws.Range("M107").Value = -1421.3079
$ws.Range("N107").Value = -60009240

$ws.Range("H123").Value = 45000
$ws.Range("J123").Value = 45000
$ws.Range("L123").Value = 45000
$ws.Range("N123").Value = -54800

$ws.Range("H137").Value = 38558
$ws.Range("J137").Value = 38558
$ws.Range("L137").Value = 38558
$ws.Range("N137").Value = -48758

$ws.Range("H139").Value = 47960
$ws.Range("J139").Value = 47960
$ws.Range("L139").Value = 47960
$ws.Range("N139").Value = -58240

$ws.Range("H141").Value = 39200
$ws.Range("J141").Value = 39200
$ws.Range("L141").Value = 39200
$ws.Range("N141").Value = -49560
